$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (row 3 through row 9) with new evaluation numbers ---

# Row 3 (Q1)
$ws.Range("B3").Value = 0.08146426224612845
$ws.Range("C3").Value = 0.4913395475084953
$ws.Range("D3").Value = 0.5384525556271494
$ws.Range("E3").Value = 0.7337932649099128
$ws.Range("F3").Value = 0.7456470948380217
$ws.Range("G3").Value = 23

# Row 4 (Q2)
$ws.Range("B4").Value = 0.5801506038144637
$ws.Range("C4").Value = 0.8159053342639054
$ws.Range("D4").Value = 4.070369851314529
$ws.Range("E4").Value = 2.017515762345992
$ws.Range("F4").Value = 1.975731026375135
$ws.Range("G4").Value = 23

# Row 5 (Q3)
$ws.Range("B5").Value = 0.2635203828962711
$ws.Range("C5").Value = 1.232511031201218
$ws.Range("D5").Value = 7.230600558527729
$ws.Range("E5").Value = 2.688977604690625
$ws.Range("F5").Value = 2.736177027977477
$ws.Range("G5").Value = 23

# Row 6 (Q4)
$ws.Range("B6").Value = 0.1794021269052618
$ws.Range("C6").Value = 1.155750348923658
$ws.Range("D6").Value = 7.068805655208574
$ws.Range("E6").Value = 2.658722560781507
$ws.Range("F6").Value = 2.712280759783512
$ws.Range("G6").Value = 23

# Row 7 (Q5)
$ws.Range("B7").Value = 0.3942472833595512
$ws.Range("C7").Value = 1.646994786995847
$ws.Range("D7").Value = 10.9495783436642
$ws.Range("E7").Value = 3.309014708892089
$ws.Range("F7").Value = 3.40075843797503
$ws.Range("G7").Value = 15

# Row 8 (Q6)
$ws.Range("B8").Value = 0.1885059718488306
$ws.Range("C8").Value = 1.762213756878201
$ws.Range("D8").Value = 11.06808830893117
$ws.Range("E8").Value = 3.326873653887561
$ws.Range("F8").Value = 3.438108980993429
$ws.Range("G8").Value = 15

# Row 9 (Q7) - now also gets an F9 value it didn't have before
$ws.Range("B9").Value = -0.1155901121422573
$ws.Range("C9").Value = 3.549470608441434
$ws.Range("D9").Value = 26.58323663793598
$ws.Range("E9").Value = 5.15589338892262
$ws.Range("F9").Value = 5.646578670017191
$ws.Range("G9").Value = 6

# --- Add new rows 10 (Q8) and 11 (Q9) ---

# Row 10 (Q8)
$ws.Range("A10").Value = "Q8"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = -2.244696353922198
$ws.Range("C10").Value = 2.375956227084492
$ws.Range("D10").Value = 17.67984107253042
$ws.Range("E10").Value = 4.204740309761165
$ws.Range("F10").Value = 4.105472664824975
$ws.Range("G10").Value = 4

# Row 11 (Q9)
$ws.Range("A11").Value = "Q9"
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").Value = -0.4278857427398495
$ws.Range("C11").Value = 0.4278857427398495
$ws.Range("D11").Value = 0.2403274097267022
$ws.Range("E11").Value = 0.4902319958210625
$ws.Range("F11").Value = 0.3383524815533931
$ws.Range("G11").Value = 2
